$wb = $excel.ActiveWorkbook

# --- "Web API (.net Framework)" sheet: mark Week 10 demo rows as Done and
# record completion dates (new column G), mirroring the pattern already
# used on earlier weeks' sheets (e.g. MYSQL!G2:G7).
$ws = $wb.Worksheets.Item("Web API (.net Framework)")
$src = $wb.Worksheets.Item("MYSQL")

# Update status column F -> "Done" for rows 2-7
$ws.Range("F2:F7").Value = "Done"

# Copy the existing date-cell number format (style index reused, no new
# style created) from MYSQL!G2 onto the target cells, then fill in the
# completion dates.
$src.Range("G2").Copy()
$ws.Range("G2:G7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G2").Value = 46027
$ws.Range("G3").Value = 46028
$ws.Range("G4").Value = 46029
$ws.Range("G5").Value = 46030
$ws.Range("G6").Value = 46030
$ws.Range("G7").Value = 46030

# --- cursor/tab bookkeeping matching the saved UI state ---
# jQuery sheet cursor moved from F23 to D15 (no data change)
$null = $wb.Worksheets.Item("jQuery").Range("D15").Select()

# The "Other" sheet was previously the active tab; the user ended on the
# "Web API (.net Framework)" sheet instead, with the cursor on F9.
$ws.Activate()
$null = $ws.Range("F9").Select()
